$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = @'
e000
'@
$ws.Range('B1').Value = @'
<Bold>e000 Welcome to Patton's Best Solo Tank Battle Game</Bold>
<LineBreak/><LineBreak/>
The game starts with a tutorial how to play. However, before starting, it is important to know that Active events are shown with a green background. The game may only advance when a green background is displayed. Most often, the game advances by rolling dice or clicking an image. 
<LineBreak/><LineBreak/>
You can explore what may happen by showing inactive events. Inactive events have a gray background. To return to the current active event, select the active event button in the status bar per the image.
<LineBreak/>
               <InlineUIContainer><Image Name='Tutorial0' Height='70'  Width='370'> </Image></InlineUIContainer>
<LineBreak/><LineBreak/>
<InlineUIContainer><Button Name='Read_Rules' Content='Read Rules' FontFamily='Courier New'  FontSize='12'> </Button></InlineUIContainer> or <InlineUIContainer><Button Name='Begin' Content='Begin Game' FontFamily='Courier New'  FontSize='12'> </Button></InlineUIContainer>
'@
$ws.Rows.Item(1).RowHeight = 180

$ws.Range('A2').Value = @'
e001
'@
$ws.Range('B2').Value = @'
<Bold>e001 Fourth Armor Division Campaign</Bold> 
<InlineUIContainer><Button Content='r1.1' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/>
The campaign game of <Bold>Patton' Best</Bold> recreates the actions of the 4th Armored Division from late July 1944 through April 1945. 
<LineBreak/><LineBreak/>
Each day begins with a check of the Combat 
<InlineUIContainer><Button Content='Calendar' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> to see what the Division was doing on that day. The four possibilities are Refitting 
<InlineUIContainer><Button Content='r27.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>, an Advance scenario <InlineUIContainer><Button Content='r20.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>, a Battle scenario 
<InlineUIContainer><Button Content='r20.3' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>, or a Counterattack scenario 
<InlineUIContainer><Button Content='r20.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.  Click image to continue.
<LineBreak/><LineBreak/>
                                            <InlineUIContainer><Image Name='Continue001' Height='100' Width='100'></Image></InlineUIContainer>
'@
$ws.Rows.Item(2).RowHeight = 210

$ws.Range('A3').Value = @'
e002
'@
$ws.Range('B3').Value = @'
<Bold>e002 Movement Board</Bold> <InlineUIContainer><Button Content='r2.11' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/>
The movement board is a depiction fo typical European countryside and is used to show the "big picture" for the day. The movement board is divided into white lines into areas. Click image to continue.
<LineBreak/><LineBreak/>
A=Farms    B=Fields    C=Villiages  D=Woods<LineBreak/>
#=Starting or exiting areas
<LineBreak/><LineBreak/>
                                   <InlineUIContainer><Image Name='MapMovement'  Height='200' Width='200'></Image></InlineUIContainer>
'@
$ws.Rows.Item(3).RowHeight = 120.75

$ws.Range('A4').Value = @'
e003
'@
$ws.Range('B4').Value = @'
<Bold>e003 Battle Board</Bold> <InlineUIContainer><Button Content='r2.12' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/>
The battle board is an abstract display used to resolve engagements with enemy forces. Your tank is placed in the center of this display and the action of an engagement revolves around it through the use of pieces representing enemy units and other informational markers. A detailed explanation is given in <InlineUIContainer><Button Content='r5.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>. Click image to continue.
<LineBreak/><LineBreak/>
                                   <InlineUIContainer><Image Name='MapBattle'  Height='200' Width='200'></Image></InlineUIContainer>
'@
$ws.Rows.Item(4).RowHeight = 105

$ws.Range('A5').Value = @'
e004
'@
$ws.Range('B5').Value = @'
<Bold>e004 Tank Card</Bold> <InlineUIContainer><Button Content='r2.13' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/>
The upper right image is the Tank Card. The game starts with the basic M4 Sherman tank, i.e., Tank Card #1. 
The Tank Card shows the tank model and other important information regarding the tank. The use of the Tank Card is described in 
<InlineUIContainer><Button Content='r5.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>. Click image to continue.
<LineBreak/><LineBreak/>
                                 <InlineUIContainer><Image Name='m001M4'  Height='200' Width='200'></Image></InlineUIContainer>
'@
$ws.Rows.Item(5).RowHeight = 105

$ws.Range('A6').Value = @'
e005
'@
$ws.Range('B6').Value = @'
<Bold>e005 After Action Report (AAR)</Bold> <InlineUIContainer><Button Content='r2.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/>
The events of each engagement or day of battle are recorded as they unfold on the After Action Report. At this time, you may elect to change the name of the tank or the names of your crew by clicking on the appropriate location on the form. 
<LineBreak/><LineBreak/>When ready, click image below to assign crew ratings to your new crew per 
<InlineUIContainer><Button Content='r7.1' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/>
                                            <InlineUIContainer><Image Name='Continue005' Height='100' Width='100'></Image></InlineUIContainer>
'@
$ws.Rows.Item(6).RowHeight = 120

$ws.Range('A7').Value = @'
e006
'@
$ws.Range('B7').Value = @'
<Bold>e006 Combat Calendar Check</Bold> 
<InlineUIContainer><Button Content='r4.1' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
<InlineUIContainer><Button Content='Calendar' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/>
Roll for possible combat today. If die &lt;= probability, start morning briefing per 
<InlineUIContainer><Button Content='e007' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.  
Otherwise continue with next day check.
<LineBreak/><LineBreak/>
Date from Combat Calendar: DATE<LineBreak/>
Expected Resistance: RESISTANCE<LineBreak/>
Probablility of Combat: PROBABILITY &gt;= 
<InlineUIContainer><Image Name='DieRollWhite' Height='21' Width='21' > </Image></InlineUIContainer>
<LineBreak/>
'@
$ws.Rows.Item(7).RowHeight = 195

$ws.Range('A8').Value = @'
e007
'@
$ws.Range('B8').Value = @'
<Bold>e007 Morning Briefing - Weather Roll</Bold> <InlineUIContainer><Button Content='r4.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/>
The 
<InlineUIContainer><Button Content='Weather' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
 Table determines weather for today:  
<InlineUIContainer><Image Name='DieRollBlue' Height='21' Width='21' > </Image></InlineUIContainer>
<LineBreak/><LineBreak/>
'@
$ws.Rows.Item(8).RowHeight = 105

$ws.Range('A9').Value = @'
e008
'@
$ws.Range('B9').Value = @'
<Bold>e008 Type of Snow</Bold> 
<LineBreak/><LineBreak/>
Snow is in the forecast. Roll for type of snow on the 
<InlineUIContainer><Button Content='Weather' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
 Table:  
<InlineUIContainer><Image Name='DieRollWhite' Height='21' Width='21' > </Image></InlineUIContainer>
<LineBreak/><LineBreak/>
'@
$ws.Rows.Item(9).RowHeight = 105

$ws.Range('A10').Value = @'
e009
'@
$ws.Range('B10').Value = @'
<Bold>e009 Ammo Loading Limits</Bold> <InlineUIContainer><Button Content='r16.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/>
See 
<InlineUIContainer><Button Content='r16.1' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
for ammo types. See 
<InlineUIContainer><Button Content='r16.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
for loading ammo. The Tank Card limits the number of normal main gun ammo allowed to AMMO_NORMAL_LOAD. Extra ammo is added in a later step after assigning normal load.
<LineBreak/><LineBreak/>
<Bold>-- AP:</Bold> Unlimited<LineBreak/>
<Bold>-- HE:</Bold> Unlimited
'@
$ws.Rows.Item(10).RowHeight = 150

$ws.Range('A11').Value = @'
e010
'@
$ws.Range('B11').Value = @'
<Bold>e010 Time Check</Bold> 
<InlineUIContainer><Button Content='r4.3' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<InlineUIContainer><Button Content='r21.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
Determine sunrise and sunset for current month using the <InlineUIContainer><Button Content='Time' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Table. 
Roll 1D on the 
<InlineUIContainer><Button Content='Time' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  Table. 
The Time Table also provides the timed used for each action take. Additionally, the same die roll is used to determine the ammo expended:  
<InlineUIContainer><Image Name='DieRollWhite' Height='21' Width='21' > </Image></InlineUIContainer>
<LineBreak/><LineBreak/>
'@
$ws.Rows.Item(11).RowHeight = 150

$ws.Range('A12').Value = @'
e011
'@
$ws.Range('B12').Value = @'
<Bold>e011 Deployment</Bold> 
<InlineUIContainer><Button Content='r4.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
Determine your tank&apos;s deployment from the 
<InlineUIContainer><Button Content='Deployment' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  Table:  
<InlineUIContainer><Image Name='DieRollBlue' Height='21' Width='21' > </Image></InlineUIContainer>
<LineBreak/><LineBreak/>
'@
$ws.Rows.Item(12).RowHeight = 105

$ws.Range('A13').Value = @'
e012
'@
$ws.Range('B13').Value = @'
<Bold>e012 Hatches</Bold> 
<InlineUIContainer><Button Content='r4.42' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
Left click on hatches on the Tank Card to toggle adding counter. If already open, click button to remove open hatch. 
<LineBreak/><LineBreak/>
                                                  <InlineUIContainer><Image Name='c15OpenHatch'  Height='80' Width='80'></Image></InlineUIContainer>
<LineBreak/><LineBreak/>
Click image in this text box when satisfied and want to continue.
'@
$ws.Rows.Item(13).RowHeight = 120

$ws.Range('A14').Value = @'
e013
'@
$ws.Range('B14').Value = @'
<Bold>e013 Gun Load</Bold> 
<InlineUIContainer><Button Content='r4.43' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
Mark the type of round you want loaded in the main gun before any action begins by clicking the highlighted box on the Tank Card in the correct ammo type box. Click image below to continue.
<LineBreak/><LineBreak/>
                                                  <InlineUIContainer><Image Name='c17GunLoad'  Height='80' Width='80'></Image></InlineUIContainer>
'@
$ws.Rows.Item(14).RowHeight = 90

$ws.Range('A15').Value = @'
e014
'@
$ws.Range('B15').Value = @'
<Bold>e014 Tank &amp; Turret Orientation</Bold> 
<InlineUIContainer><Button Content='r4.44' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
The Sherman tank counter is placed at the center of the battle board. If you want the turret to face a different sector, click tank counter on center of Battle Board. 
Alternatively, select buttons here:
<LineBreak/><LineBreak/>
                                   <InlineUIContainer><Button Content='  -  ' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
<InlineUIContainer><Image Name='c16Turret'  Height='150' Width='150'></Image></InlineUIContainer> 
<InlineUIContainer><Button Content='  +  ' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/>
When you are satisfied with the current turret orientation, click turret image between buttons to continue.
'@
$ws.Rows.Item(15).RowHeight = 165

$ws.Range('A16').Value = @'
e015
'@
$ws.Range('B16').Value = @'
<Bold>e015 Loader Spotting</Bold> 
<InlineUIContainer><Button Content='r4.45' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<InlineUIContainer><Button Content='r17.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
If the loader is buttoned up or does not have a hatch, mark the sector he will be searching by left clicking on dot just outside the long range zone of the sector. 
<LineBreak/><LineBreak/>
'@
$ws.Rows.Item(16).RowHeight = 90

$ws.Range('A17').Value = @'
e016
'@
$ws.Range('B17').Value = @'
<Bold>e016 Commander Spotting</Bold> 
<InlineUIContainer><Button Content='r4.45' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<InlineUIContainer><Button Content='r17.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
If the Commander is buttoned up or does not have a vision cupola, mark the sector he will be searching by right clicking on dot just outside the long range zone of the sector. 
<LineBreak/><LineBreak/>
'@
$ws.Rows.Item(17).RowHeight = 90

$ws.Range('A18').Value = @'
e017
'@
$ws.Range('B18').Value = @'
<Bold>e017 Preparations Final</Bold> 
<InlineUIContainer><Button Content='r4.46' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
US Control markers are placed on sectors 1, 2, and 3. The Weather is displayed on top left of Battle Board.
<LineBreak/><LineBreak/>
                                            <InlineUIContainer><Image Name='Continue017' Height='100' Width='100'></Image></InlineUIContainer>
'@
$ws.Rows.Item(18).RowHeight = 90

$ws.Range('A19').Value = @'
e018
'@
$ws.Range('B19').Value = @'
<Bold>e018 Set Start Area</Bold> 
<InlineUIContainer><Button Content='r4.51' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Around the edge of the Movement Board, there are 10 areas number 1-10. The area is marked with the Start Area and Task Force markers. 
<LineBreak/><LineBreak/>
Die Roll =  <InlineUIContainer><Image Name='DieRollWhite' Height='21' Width='21' > </Image></InlineUIContainer>
<LineBreak/><LineBreak/>
'@
$ws.Rows.Item(19).RowHeight = 105

$ws.Range('A20').Value = @'
e019
'@
$ws.Range('B20').Value = @'
<Bold>e019 Set Exit Area</Bold> 
<InlineUIContainer><Button Content='r4.52' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
On the <InlineUIContainer><Button Content='Exit Areas' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
Table, roll 1D and cross reference the number with the Start Area marker 
<InlineUIContainer><Button Content='r4.51' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.  
<LineBreak/><LineBreak/>
Die Roll =  <InlineUIContainer><Image Name='DieRollWhite' Height='21' Width='21' > </Image></InlineUIContainer> 
<LineBreak/><LineBreak/>
'@
$ws.Rows.Item(20).RowHeight = 135

$ws.Range('A21').Value = @'
e020
'@
$ws.Range('B21').Value = @'
<Bold>e020 Enemy Strength Check - Choose Area</Bold> 
<InlineUIContainer><Button Content='r4.53' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Check any one adjacent area to your task force for estimating enemy strength. Click on one of the adjacent regions highlighted blue.
<LineBreak/><LineBreak/>
'@
$ws.Rows.Item(21).RowHeight = 75

$ws.Range('A22').Value = @'
e021
'@
$ws.Range('B22').Value = @'
<Bold>e021 Enemy Strength Check Roll</Bold> 
<InlineUIContainer><Button Content='r4.53' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Roll 1D and consult the <InlineUIContainer><Button Content='Resistance' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
Table. The area is marked with a Light, Medium, or Heavy marker.
<LineBreak/><LineBreak/>
Die Roll =  <InlineUIContainer><Image Name='DieRollWhite' Height='21' Width='21' > </Image></InlineUIContainer> 
<LineBreak/><LineBreak/>
'@
$ws.Rows.Item(22).RowHeight = 120

$ws.Range('A23').Value = @'
e022
'@
$ws.Range('B23').Value = @'
<Bold>e022 Choose Operations</Bold> 
<InlineUIContainer><Button Content='r4.54' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Choose one of following options. To see the options, select the &apos;e###&apos; button. To choose the option, select the other buttons. Each option uses up time per the 
<InlineUIContainer><Button Content='Time' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
Table.
<LineBreak/><LineBreak/>
<InlineUIContainer><Button Content='e020' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
 Perform an additional Enemy <InlineUIContainer><Button Content='Strength Check' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer><LineBreak/>
<InlineUIContainer><Button Content='e023' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Call for Artillery 
<InlineUIContainer><Button Content='Support' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer><LineBreak/>
<InlineUIContainer><Button Content='e025' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Call for Air 
<InlineUIContainer><Button Content='Strike' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer><LineBreak/>
<InlineUIContainer><Button Content='e027' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Attempt to 
<InlineUIContainer><Button Content='Resupply' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer><LineBreak/>
<InlineUIContainer><Button Content='e028' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Choose adjacent area to 
<InlineUIContainer><Button Content='Enter' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer><LineBreak/>
'@
$ws.Rows.Item(23).RowHeight = 255

$ws.Range('A24').Value = @'
e023
'@
$ws.Range('B24').Value = @'
<Bold>e023 Call for Artillery Support</Bold> 
<InlineUIContainer><Button Content='r23.1' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<InlineUIContainer><Button Content='r4.54.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Call to hit an area adjacent to your task force. Only three Artillery Support markers can exist on the board at one time. Consult the 
<InlineUIContainer><Button Content='Time' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
to see if Artillery Support arrives. If successful, an Artillery Support marker on the area. Click on one of the adjacent regions highlighted blue.
'@
$ws.Rows.Item(24).RowHeight = 105

$ws.Range('A25').Value = @'
e024
'@
$ws.Range('B25').Value = @'
<Bold>e024 Artillery Support Roll</Bold> 
<InlineUIContainer><Button Content='r23.1' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Roll 1D and consult the <InlineUIContainer><Button Content='Time' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
Table to see if Artillery Support arrives. If successful, an Artillery Support marker is placed on the area. 
<LineBreak/><LineBreak/>
Die Roll =  <InlineUIContainer><Image Name='DieRollWhite' Height='21' Width='21' > </Image></InlineUIContainer> 
<LineBreak/><LineBreak/>
'@
$ws.Rows.Item(25).RowHeight = 120

$ws.Range('A26').Value = @'
e025
'@
$ws.Range('B26').Value = @'
<Bold>e025 Call for Air Strike</Bold> 
<InlineUIContainer><Button Content='r23.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<InlineUIContainer><Button Content='r4.54.3' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Call to hit an area adjacent to your task force. Only two Air Strike markers can exist on the board at one time. 
An additional 15 minute action can be selected while waiting for an air strike. Consult the 
<InlineUIContainer><Button Content='Time' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
to see if Air Strike arrives. If successful, an Air Strike marker is placed on the area. Click on one of the adjacent regions highlighted blue.
'@
$ws.Rows.Item(26).RowHeight = 120

$ws.Range('A27').Value = @'
e026
'@
$ws.Range('B27').Value = @'
<Bold>e026 Air Strike Roll</Bold> 
<InlineUIContainer><Button Content='r4.54.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Roll 1D and consult the <InlineUIContainer><Button Content='Time' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
Table to see if Air Strike arrives. If successful, an Air Strike marker is placed on the area. 
<LineBreak/><LineBreak/>
Die Roll =  <InlineUIContainer><Image Name='DieRollWhite' Height='21' Width='21' > </Image></InlineUIContainer> 
<LineBreak/><LineBreak/>
'@
$ws.Rows.Item(27).RowHeight = 120

$ws.Range('A28').Value = @'
e027
'@
$ws.Range('B28').Value = @'
<Bold>e027 Attempt to Resupply</Bold> 
<InlineUIContainer><Button Content='r4.54.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Consult the 
<InlineUIContainer><Button Content='Time' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
to see if resupply occurs. If successful, you may relead your tank with ammo. 
<LineBreak/><LineBreak/>
Die Roll =  <InlineUIContainer><Image Name='DieRollWhite' Height='21' Width='21' > </Image></InlineUIContainer> 
<LineBreak/><LineBreak/>
'@
$ws.Rows.Item(28).RowHeight = 135

$ws.Range('A29').Value = @'
e028
'@
$ws.Range('B29').Value = @'
<Bold>e028 Enter Adjacent Area</Bold> 
<InlineUIContainer><Button Content='r4.54.5' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Click on one of the adjacent highlighted areas. Artillery Support or Air Strike Counters are moved to the battle board as a reminder. 
<LineBreak/><LineBreak/>
                        <InlineUIContainer><Image Name='Sherman1' Height='200' Width='325'></Image></InlineUIContainer>
'@
$ws.Rows.Item(29).RowHeight = 90

$ws.Range('A30').Value = @'
e029
'@
$ws.Range('B30').Value = @'
<Bold>e029 Advancing Fire Choice</Bold> 
<InlineUIContainer><Button Content='r4.54.5' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Choose one of the following images for Advancing Fire option per 
<InlineUIContainer><Button Content='r22.1' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>. 
<LineBreak/><LineBreak/>
<InlineUIContainer><Image Name='c44AdvanceFire' Height='60' Width='60'></Image></InlineUIContainer> to perform advancing fire. <LineBreak/> <LineBreak/>
<InlineUIContainer><Image Name='c44AdvanceFireDeny' Height='60' Width='60'></Image></InlineUIContainer> to skip advancing fire.  
'@
$ws.Rows.Item(30).RowHeight = 120

$ws.Range('A31').Value = @'
e030
'@
$ws.Range('B31').Value = @'
<Bold>e030 Advancing Fire Ammo Use</Bold> 
<InlineUIContainer><Button Content='r22.11' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Mark off 1D/2 (round down) HE rounds and .30 caliber MG ammo boxes regardless of whether the battle occurs or not. Mark off on the After Action Report 
<InlineUIContainer><Button Content='AAR' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.
<LineBreak/><LineBreak/>
Die Roll =  <InlineUIContainer><Image Name='DieRollWhite' Height='21' Width='21' > </Image></InlineUIContainer> 
'@
$ws.Rows.Item(31).RowHeight = 105

$ws.Range('A32').Value = @'
e031
'@
$ws.Range('B32').Value = @'
<Bold>e031 Enemy Strength Roll Entering Battle Board</Bold> 
<InlineUIContainer><Button Content='r4.53' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Roll 1D and consult the <InlineUIContainer><Button Content='Resistance' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
Table for enemy strength: <LineBreak/><LineBreak/>
Die Roll =  <InlineUIContainer><Image Name='DieRollWhite' Height='21' Width='21' > </Image></InlineUIContainer> 
'@
$ws.Rows.Item(32).RowHeight = 90

$ws.Range('A33').Value = @'
e032
'@
$ws.Range('B33').Value = @'
<Bold>e032 Battle Check</Bold> 
<InlineUIContainer><Button Content='r4.54.5' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Roll 1D and consult the <InlineUIContainer><Button Content='Resistance' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
Table to determine if combat occurs in this area: <LineBreak/><LineBreak/>
Die Roll =  <InlineUIContainer><Image Name='DieRollWhite' Height='21' Width='21' > </Image></InlineUIContainer> 
'@
$ws.Rows.Item(33).RowHeight = 90

$ws.Range('A34').Value = @'
e033
'@
$ws.Range('B34').Value = @'
<Bold>e033 No Combat</Bold> 
<InlineUIContainer><Button Content='r4.54.5' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
If converting territory to US Control, Victory points are added to the After Action Report 
<InlineUIContainer><Button Content='AAR' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.
'@
$ws.Rows.Item(34).RowHeight = 75

$ws.Range('A35').Value = @'
e034
'@
$ws.Range('B35').Value = @'
<Bold>e034 Placing Advancing Fire Markers</Bold> 
<InlineUIContainer><Button Content='r4.61' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Place Advancing Fire Markers available to you per 
<InlineUIContainer><Button Content='r22.12' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>. 
Click one of highlighted regions to place. 
<LineBreak/><LineBreak/>
You place up to six minus one marker for every three friendly tank losses (rounded up) . You may place more than one in a zone. 
<LineBreak/><LineBreak/>
The status bar on the bottom shows how many are remaining to place.
'@
$ws.Rows.Item(35).RowHeight = 150

$ws.Range('A36').Value = @'
e035
'@
$ws.Range('B36').Value = @'
<Bold>e035 Ambush Check</Bold> 
<InlineUIContainer><Button Content='r4.65' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Roll 1D for possible Ambush: 
<InlineUIContainer><Image Name='DieRollWhite' Height='21' Width='21' > </Image></InlineUIContainer>&lt; 8
<LineBreak/><LineBreak/>
'@
$ws.Rows.Item(36).RowHeight = 90

$ws.Range('A37').Value = @'
e036
'@
$ws.Range('B37').Value = @'
<Bold>e036 Battle Board Empty</Bold> 
<InlineUIContainer><Button Content='r4.77' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Since the Battle Board is now empty of enemy units, the battle for this area is over. 
<LineBreak/><LineBreak/>
1.) Flip Resistance marker to US Controlled on Movement Board.
<LineBreak/>
2.) Victory points for control of the area added to the 
<InlineUIContainer><Button Content='AAR' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.
<LineBreak/>
3.) If daylight remains, return to Prepare for Battle per 
<InlineUIContainer><Button Content='r4.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.
<LineBreak/>
4.) No daylight, perform the Evening Debriefing per 
<InlineUIContainer><Button Content='r4.9' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.
<LineBreak/>
5.) Click image to continue.
<LineBreak/><LineBreak/>
          <InlineUIContainer><Image Name='Debrief' Height='225' Width='450'></Image></InlineUIContainer>
'@
$ws.Rows.Item(37).RowHeight = 285

$ws.Range('A38').Value = @'
e037
'@
$ws.Range('B38').Value = @'
<Bold>e037 Smoke Depletion Phase</Bold> 
<InlineUIContainer><Button Content='r4.71' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Deplete smoke in each zone by converting one white full strength Smoke marker to a gray 1/2 strength Smoke marker. Alternatively, remove 1/2 strength Smoke marker. Refer to 
<InlineUIContainer><Button Content='r18.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> for the smoke rules. Click image to continue with 
<InlineUIContainer><Button Content='r4.72' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.
<LineBreak/><LineBreak/>
                                              <InlineUIContainer><Image Name='c111Smoke1' Height='100' Width='100'></Image></InlineUIContainer>
'@
$ws.Rows.Item(38).RowHeight = 120

$ws.Range('A39').Value = @'
e038
'@
$ws.Range('B39').Value = @'
<Bold>e038 Orders Phase</Bold> 
<InlineUIContainer><Button Content='r4.73' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Click on blue squares to open hatches. Click on Open Hatch marker to close.
 <LineBreak/><LineBreak/>
Click crew member action boxes to select from a pull down to assign crew actions for each crew member per 
<InlineUIContainer><Button Content='r8.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.
<LineBreak/><LineBreak/>
Click on the appropriate Gun Load box to set the Gun Reload marker per 
<InlineUIContainer><Button Content='r5.23' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.
<LineBreak/><LineBreak/>
Click the Gun Reload marker/button if you want the reload to come from the ready rack. A Ready Rack Ammo Reload marker is added per  
<InlineUIContainer><Button Content='r9.6' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.
<LineBreak/><LineBreak/>
Determine the specific unit type for any units identified per 
<InlineUIContainer><Button Content='r17.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.
<LineBreak/><LineBreak/>
'@
$ws.Rows.Item(39).RowHeight = 255

$ws.Range('A40').Value = @'
e039
'@
$ws.Range('B40').Value = @'
<Bold>e039 Random Events for Ambush</Bold> 
<InlineUIContainer><Button Content='r4.65' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
<LineBreak/><LineBreak/>
Roll on the 
<InlineUIContainer><Button Content='Random Events' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
Table:  
<InlineUIContainer><Image Name='DieRollBlue' Height='21' Width='21' > </Image></InlineUIContainer>
<LineBreak/><LineBreak/>
  <InlineUIContainer><Button Content='e040' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Time Passes<LineBreak/>
  <InlineUIContainer><Button Content='e041' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Friendly Artillery<LineBreak/>
  <InlineUIContainer><Button Content='e042' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Enemy Artillery<LineBreak/>
  <InlineUIContainer><Button Content='e043' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Mine Attack<LineBreak/>
  <InlineUIContainer><Button Content='e044' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Panzerfaust Attack<LineBreak/>
  <InlineUIContainer><Button Content='e045' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Harrassing Fire<LineBreak/>
  <InlineUIContainer><Button Content='e046' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Friendly Advance<LineBreak/>
  <InlineUIContainer><Button Content='e047' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Enemy Reinforcment<LineBreak/>
  <InlineUIContainer><Button Content='e048' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Enemy Advance<LineBreak/>
  <InlineUIContainer><Button Content='e049' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Flanking Fire<LineBreak/>
<LineBreak/><LineBreak/>
'@
$ws.Rows.Item(40).RowHeight = 285

$ws.Range('A41').Value = @'
e040
'@
$ws.Range('B41').Value = @'
<Bold>e040 Time Passes</Bold> 
<InlineUIContainer><Button Content='r21.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
Fifteen minutes pass on the After Action Report 
<InlineUIContainer><Button Content='AAR' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>. 
<LineBreak/><LineBreak/>
Click image to continue.
<LineBreak/><LineBreak/>
                                   <InlineUIContainer><Image Name='MilitaryWatch' Height='100' Width='200'></Image></InlineUIContainer>
'@
$ws.Rows.Item(41).RowHeight = 135

$ws.Range('A42').Value = @'
e041
'@
$ws.Range('B42').Value = @'
<Bold>e041 Friendly Artillery</Bold> 
<InlineUIContainer><Button Content='r23.1' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
Friendly artillery support arrives. Click image to continue.
<LineBreak/><LineBreak/>
                                        <InlineUIContainer><Image Name='c39ArtillerySupport'  Height='80' Width='80'></Image></InlineUIContainer>
'@
$ws.Rows.Item(42).RowHeight = 90

$ws.Range('A43').Value = @'
e042
'@
$ws.Range('B43').Value = @'
<Bold>e042 Enemy Artillery</Bold> 
<LineBreak/><LineBreak/>
Enemy artillery arrives. Roll 1D to knock out (KO) friendly units: 
<InlineUIContainer><Image Name='DieRollWhite' Height='21' Width='21' > </Image></InlineUIContainer>
<LineBreak/><LineBreak/>
1 KO for 1-6<LineBreak/>
2 KO for 7-9<LineBreak/>
3 KO for 10
<LineBreak/><LineBreak/>
'@
$ws.Rows.Item(43).RowHeight = 135

$ws.Range('A44').Value = @'
e043
'@
$ws.Range('B44').Value = @'
<Bold>e043 Mine Attack</Bold> 
<InlineUIContainer><Button Content='r15.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
Roll 1D on the 
<InlineUIContainer><Button Content='Minefield' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Attack Table: 
<InlineUIContainer><Image Name='DieRollWhite' Height='21' Width='21' > </Image></InlineUIContainer> 
'@
$ws.Rows.Item(44).RowHeight = 90

$ws.Range('A45').Value = @'
e043a
'@
$ws.Range('B45').Value = @'
<Bold>e043a Mine Attack Ignored</Bold> 
<InlineUIContainer><Button Content='r15.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
No effect since Sherman is not moving. Click image to continue.
<LineBreak/><LineBreak/>
                                            <InlineUIContainer><Image Name='Continue043' Height='100' Width='100'></Image></InlineUIContainer>
'@
$ws.Rows.Item(45).RowHeight = 90

$ws.Range('A46').Value = @'
e044
'@
$ws.Range('B46').Value = @'
<Bold>e044 Panzerfaust Attack</Bold> 
<InlineUIContainer><Button Content='r15.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
Determine from which sector of the Battle Board attack is originating by rolling 1D according to 
<InlineUIContainer><Button Content='r5.12' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> :  
<InlineUIContainer><Image Name='DieRollWhite' Height='21' Width='21' > </Image></InlineUIContainer>
<LineBreak/><LineBreak/>
If the sector rolled is US controlled, no attack is made. If the sector is not US Controlled, a Panzerfaust marker is placed in the sector's close range.
<LineBreak/><LineBreak/>
If an attack occurs, roll on the 
<InlineUIContainer><Button Content='Panzerfaust' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  Attack Table to hit Sherman:  
<InlineUIContainer><Image Name='DieRollWhite1' Height='21' Width='21' Visibility='Hidden'> </Image></InlineUIContainer>
<LineBreak/><LineBreak/>
If hits, roll again to see if your tank is knocked out (KO):  
<InlineUIContainer><Image Name='DieRollWhite2' Height='21' Width='21' Visibility='Hidden'> </Image></InlineUIContainer>
<LineBreak/><LineBreak/>
'@
$ws.Rows.Item(46).RowHeight = 240

$ws.Range('A47').Value = @'
e045
'@
$ws.Range('B47').Value = @'
<Bold>e045 Harrassing Fire</Bold> 
<InlineUIContainer><Button Content='r15.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
Your tank is sprayed with small weapons fire. Roll 1D on the 
<InlineUIContainer><Button Content='Collateral' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Damage Table:  
<InlineUIContainer><Image Name='DieRollWhite' Height='21' Width='21' > </Image></InlineUIContainer>
<LineBreak/><LineBreak/>
'@
$ws.Rows.Item(47).RowHeight = 105

$ws.Range('A48').Value = @'
e046
'@
$ws.Range('B48').Value = @'
<Bold>e046 Friendly Advance</Bold> 
<InlineUIContainer><Button Content='r15.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
Place a US Controlled marker in a sector of your choice that is empty of enemy units and adjacent to a sector already US Controlled.
'@
$ws.Rows.Item(48).RowHeight = 60

$ws.Range('A49').Value = @'
e046a
'@
$ws.Range('B49').Value = @'
<Bold>e046a Friendly Advance Ignored</Bold> 
<InlineUIContainer><Button Content='r15.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
Since there is no sector adjacent to a US Controlled sector that is empty of enemy units, this event is ignored. Click image to continue.
'@
$ws.Rows.Item(49).RowHeight = 60

$ws.Range('A50').Value = @'
e047
'@
$ws.Range('B50').Value = @'
<Bold>e047 Enemy Reinformement</Bold> 
<InlineUIContainer><Button Content='r15.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
'@
$ws.Rows.Item(50).RowHeight = 45

$ws.Range('A51').Value = @'
e048
'@
$ws.Range('B51').Value = @'
<Bold>e048 Enemy Advance</Bold> 
<InlineUIContainer><Button Content='r15.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
Remove one US Control marker from sector adjacent to an enemy unit. If two sectors are eligible, it is chosen randomly
'@
$ws.Rows.Item(51).RowHeight = 60

$ws.Range('A52').Value = @'
e048a
'@
$ws.Range('B52').Value = @'
<Bold>e048a Enemy Advance Ignored</Bold> 
<InlineUIContainer><Button Content='r15.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
Since no sector is adjacent to an enemy unit, this event is ignored. Click image to continue.
'@
$ws.Rows.Item(52).RowHeight = 60

$ws.Range('A53').Value = @'
e049
'@
$ws.Range('B53').Value = @'
<Bold>e049 Flanking Fire</Bold> 
<InlineUIContainer><Button Content='r15.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
Roll against each enemy unit on the Friendly Action Table with a die roll modifer of -10. Click image to continue.
<LineBreak/><LineBreak/>
'@
$ws.Rows.Item(53).RowHeight = 75

$ws.Range('A54').Value = @'
e050
'@
$ws.Range('B54').Value = @'
<Bold>e050 Evening Debriefing</Bold> 
<LineBreak/><LineBreak/>
An evening debriefing is performed per <InlineUIContainer><Button Content='r4.9' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>. 
Click image to continue.
<LineBreak/><LineBreak/>
                     <InlineUIContainer><Image Name='Sherman4' Height='168' Width='275'></Image></InlineUIContainer>
'@
$ws.Rows.Item(54).RowHeight = 90

$ws.Range('A55').Value = @'
e501
'@
$ws.Range('B55').Value = @'
<Bold>e501 Game Won!!!</Bold>
<LineBreak/><LineBreak/>
'@
$ws.Rows.Item(55).RowHeight = 30

$ws.Range('A56').Value = @'
e502
'@
$ws.Range('B56').Value = @'
<Bold>e502 Game Lost</Bold>
<LineBreak/><LineBreak/>
'@
$ws.Rows.Item(56).RowHeight = 30

$ws.Range('A57').Value = @'
e503
'@
$ws.Range('B57').Value = @'
<Bold>e503 End Game Statistics and Feats </Bold>
<LineBreak/><LineBreak/>
Select 'File | New' menu option to play again.
<LineBreak/><LineBreak/>
Click image to review map: 
<LineBreak/><InlineUIContainer><Image Name='Win' Height='100' Width='100'></Image></InlineUIContainer>
<LineBreak/><LineBreak/>
Click image to exit game:
<LineBreak/>   <InlineUIContainer><Image Name='DoorClosing' Height='150' Width='75'></Image></InlineUIContainer>
'@
$ws.Rows.Item(57).RowHeight = 135

$ws.Range('A58').Value = @'
e503a
'@
$ws.Range('B58').Value = @'
<Bold>e503a Game Feats </Bold>
<LineBreak/><LineBreak/>You completed an achievement that requires great courage, skill, and persistence.
<LineBreak/><LineBreak/>Congratulations!
<LineBreak/><LineBreak/>Click star to continue.
'@
$ws.Rows.Item(58).RowHeight = 60

$ws.Range('A59').Value = @'
e504
'@
$ws.Range('B59').Value = @'
<Bold>e504 Play Again?</Bold>
<LineBreak/><LineBreak/>
Select 'File | New' menu option to play again.
<LineBreak/><LineBreak/>
Click image to exit game:
<LineBreak/>   <InlineUIContainer><Image Name='DoorClosing'  Height='150' Width='75'></Image></InlineUIContainer>
'@
$ws.Rows.Item(59).RowHeight = 90

[void]$ws.Range('B41').Select()
